$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update odds values in row 3 (Santa Fe vs Bucaramanga) as per the latest
# FlashScore refresh:
#   Odd_Over05_FT (M3):  1.13 -> 1.11
#   Odd_Under05_FT (N3): 6    -> 6.5
#   Odd_Over3_FT (U3):   4.4  -> 4.3
#   Odd_Under3_FT (V3):  1.2  -> 1.21
$ws.Range("M3").Value2 = 1.11
$ws.Range("N3").Value2 = 6.5
$ws.Range("U3").Value2 = 4.3
$ws.Range("V3").Value2 = 1.21
